$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 09:07:22.099000 to 2024-03-11 09:55:59.893000"

$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"
$ws.Range("B2").Value = 0.03483537037037037

$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 27.75552166666667

$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1385.740676991389

$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 38.746

$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 10.259

$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 25.0

$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 97.0

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 27.14685761272991

$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 51.04608042521935

$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 72.0

$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Custom mode`n90.09%`nEco mode`n8.86%`nSports mode`n0.07%"

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 5724.067349

$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -1668.450313434448

$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 0.02443286916666667

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.001763132077205007

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.326

$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 2.963

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.363

$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 25.0

$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 40.0

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 15.0

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 70.0

$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 63.0

$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 64.0

$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 60.0

$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 95.0

$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0.0

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 40.0

$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 25.0

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 15.0

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 55.0

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.526553691666667

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.0000001409247896741873

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 37.0

$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 25.37767519932857

$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 8.203944607637432

$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 6.602322003077353

$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 7.508043082948664

$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 8.354315288851588

$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 11.63799132745839

$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 8.106028815218911

$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 12.75003496992586

$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 11.31626800951182

$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0.0
